$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-26
$data = @(
    @(9, 9),
    @(6, 7),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(3, 3),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(9, 9),
    @(10, 10),
    @(10, 10),
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(4, 5),
    @(6, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
